# UserDetails.xlsx update:
#  - add a new "password" column (I) with a sample value for every data row
#  - remove the extra duplicate rows (6-12), keeping only the first 4 data rows
#  - leave selection on the last data row, matching the author's last save

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "password" header + values for the 4 data rows that remain (rows 2-5)
$ws.Range("I1").Value = "password"
$ws.Range("I2:I5").Value = 12345678

# Drop the redundant duplicate rows 6-12 (only rows 1-5 should remain)
$ws.Rows("6:12").Delete()

# Match the saved selection (whole row 5 selected)
$ws.Rows("5:5").Select()
